# Liquidity Override workbook update - generated DIF 20201030
# liquidity and asset allocation reports
#
# Adds a new data row (row 3) to Sheet1:
#   A3 = 2020-10-30 (date, same formatting as A2)
#   B3 = "6688 HK Equity"
#   C3 = "L0"
# and moves the active selection to D3, matching the
# author's saved worksheet state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date cell above first so the new date cell inherits the
# existing date number format (style) instead of creating a new one.
$ws.Range("A2").Copy($ws.Range("A3")) | Out-Null

$ws.Range("A3").Value = 44134
$ws.Range("B3").Value = "6688 HK Equity"
$ws.Range("C3").Value = "L0"

$ws.Range("D3").Select()
